# Apply the "conversation-type" IG publish update:
#  - Metadata sheet: bump Version, Date, set Publisher, replace the
#    duplicated "Contact" rows with a single "Jurisdiction" row.
#  - Elements sheet: update the root Extension row's Short/Definition
#    text to the resource-specific description.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Drop the second (duplicate) "Contact" row (row 11); everything below
# shifts up by one, matching the new A1:B20 dimension.
$meta.Rows.Item(11).Delete()

# Version / Date bump.
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value, and the old "Contact" row becomes
# "Jurisdiction" / "United States of America".
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition text now reflects this
# extension's actual purpose instead of the generic placeholder.
$elements.Range("K2").Value = "Conversation Type"
$elements.Range("L2").Value = "Type of conversation, used for ihe chat communications"
